# "break out stock.yaml completed"
# A new scrape run (04/06/2024 02:42:39) appended its rows to each
# "breakout" sheet, and the previous run's bsecode (column D) cells -
# which had been left as text - got normalized to real numbers.
#
# Helper: write a value into a cell as TEXT even when it looks numeric
# (Excel normally auto-converts a numeric-looking string typed into
# Range.Value into a real number, same as typing it on the grid - so we
# use the leading-apostrophe text-entry convention, then strip the
# resulting quote-prefixed style back to Normal so no stray formatting
# is left behind on the cell).
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "10per change": A1:H15 -> A1:H16
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("10per change")

# D15 was stored as text "542651" -> make it numeric
$ws.Range("D15").Value = 542651

# New row 16: duplicate of row 15's stock, new scrape timestamp
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "KPITTECH"
$ws.Range("C16").Value = "KPIT Technologies Ltd"
Set-TextValue $ws.Range("D16") "542651"
$ws.Range("E16").Value = -0.9
$ws.Range("F16").Value = 1444.2
$ws.Range("G16").Value = 795009
$ws.Range("H16").Value = "04/06/2024 02:42:39"

# ---------------------------------------------------------------
# Sheet "3 V 0.3": A1:H17 -> A1:H21
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("3 V 0.3")

# D14:D17 were stored as text -> make them numeric
$ws.Range("D14").Value = 532134
$ws.Range("D15").Value = 540065
$ws.Range("D16").Value = 517300
$ws.Range("D17").Value = 530965

# New row 18
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "BANKBARODA"
$ws.Range("C18").Value = "Bank Of Baroda"
Set-TextValue $ws.Range("D18") "532134"
$ws.Range("E18").Value = 12.1
$ws.Range("F18").Value = 296.95
$ws.Range("G18").Value = 74818392
$ws.Range("H18").Value = "04/06/2024 02:42:39"

# New row 19
$ws.Range("A19").Value = 2
$ws.Range("B19").Value = "RBLBANK"
$ws.Range("C19").Value = "Rbl Bank Limited"
Set-TextValue $ws.Range("D19") "540065"
$ws.Range("E19").Value = 6.33
$ws.Range("F19").Value = 261.25
$ws.Range("G19").Value = 19552284
$ws.Range("H19").Value = "04/06/2024 02:42:39"

# New row 20
$ws.Range("A20").Value = 3
$ws.Range("B20").Value = "GIPCL"
$ws.Range("C20").Value = "Gujarat Industries Power Company Limited"
Set-TextValue $ws.Range("D20") "517300"
$ws.Range("E20").Value = 16.68
$ws.Range("F20").Value = 209.15
$ws.Range("G20").Value = 7996068
$ws.Range("H20").Value = "04/06/2024 02:42:39"

# New row 21
$ws.Range("A21").Value = 4
$ws.Range("B21").Value = "IOC"
$ws.Range("C21").Value = "Indian Oil Corporation Limited"
Set-TextValue $ws.Range("D21") "530965"
$ws.Range("E21").Value = 7.94
$ws.Range("F21").Value = 175.3
$ws.Range("G21").Value = 62214295
$ws.Range("H21").Value = "04/06/2024 02:42:39"

# ---------------------------------------------------------------
# Sheet "DND 3 V 0.3": A1:H12 -> A1:H15
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("DND 3 V 0.3")

# D10:D12 were stored as text -> make them numeric
$ws.Range("D10").Value = 532134
$ws.Range("D11").Value = 517300
$ws.Range("D12").Value = 532885

# New row 13
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "BANKBARODA"
$ws.Range("C13").Value = "Bank Of Baroda"
Set-TextValue $ws.Range("D13") "532134"
$ws.Range("E13").Value = 12.1
$ws.Range("F13").Value = 296.95
$ws.Range("G13").Value = 74818392
$ws.Range("H13").Value = "04/06/2024 02:42:39"

# New row 14
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "GIPCL"
$ws.Range("C14").Value = "Gujarat Industries Power Company Limited"
Set-TextValue $ws.Range("D14") "517300"
$ws.Range("E14").Value = 16.68
$ws.Range("F14").Value = 209.15
$ws.Range("G14").Value = 7996068
$ws.Range("H14").Value = "04/06/2024 02:42:39"

# New row 15
$ws.Range("A15").Value = 3
$ws.Range("B15").Value = "CENTRALBK"
$ws.Range("C15").Value = "Central Bank Of India"
Set-TextValue $ws.Range("D15") "532885"
$ws.Range("E15").Value = 11.83
$ws.Range("F15").Value = 72.3
$ws.Range("G15").Value = 95389337
$ws.Range("H15").Value = "04/06/2024 02:42:39"
